# build: update test database with dummy step
# Populate the (until now empty) Sheet1 header/data row with the dummy
# test-run record used by the test database fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "random-id"
$ws.Range("B1").Value = "random-test-id"
$ws.Range("C1").Value = "test-action"
$ws.Range("D1").Value = "random description"
$ws.Range("E1").Value = "random_runner_name"
$ws.Range("F1").Value = "PENDING"
# Note: column G is intentionally left blank/untouched.
$ws.Range("H1").Value = "random data"

# Give the populated cells the (automatic/theme) text colour Excel applies,
# which introduces the second font + cell style used by this row.
$ws.Range("A1:F1").Font.ThemeColor = 1
$ws.Range("H1").Font.ThemeColor = 1
